$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.708.87"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.40%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.683.16"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.52%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.27%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.07%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.29%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3930"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.20%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3968"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.14%  "

# Row 9
$ws.Range("E9").Value = "  +0.31%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.425"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.51%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.69"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.85%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08697"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.30%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.34"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.53%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.347"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.98%  "

# Row 15
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001329"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.16%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.812"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.65%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.678.98"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.12%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.77"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.06%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07108"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.38%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.31"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.84%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.146"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.04%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.35%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.08"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.28%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.717.59"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.37%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.351"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.30%  "

# Row 26
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "23.59"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.44%  "

# Row 27
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.789"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -7.08%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.06"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.93%  "

# Row 29
$ws.Range("B29").Value = "HuobiToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.820"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.06%  "

# Row 30
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "150.85"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.98%  "

# Row 31
$ws.Range("E31").Value = "  -7.91%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.414"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.35%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.862.15"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.96%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08449"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.10%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.03088"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.96%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.014"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.91%  "

# Row 37
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.956"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.49%  "

# Row 38
$ws.Range("E38").Value = "  -1.73%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09540"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.23%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.57"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.81%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7972"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.59%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.476"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.07%  "

# Row 43
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.72"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.71"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.31%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7173"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.26%  "

# Row 46
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.587"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.03%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.177"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.03%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08670"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +4.47%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.003"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.24%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.335"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.33%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "138.65"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.64%  "
